$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet1: Recommandations (rows 2-48) ---
$ws1.Range("A2").Value = "BRVM - SERVICES PUBLICS"
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 8
$ws1.Range("D2").Value = 3317.77
$ws1.Range("E2").Value = 110.09
$ws1.Range("F2").Value = "🟡 Observer"
$ws1.Range("G2").Value = "➖ Neutre"

$ws1.Range("A3").Value = "CFAO MOTORS CI"
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 4
$ws1.Range("D3").Value = 2700
$ws1.Range("E3").Value = 675
$ws1.Range("F3").Value = "🟡 Observer"
$ws1.Range("G3").Value = "➖ Neutre"

$ws1.Range("A4").Value = "SAFCA CI"
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 4
$ws1.Range("D4").Value = 2690
$ws1.Range("E4").Value = 675
$ws1.Range("F4").Value = "🟡 Observer"
$ws1.Range("G4").Value = "➖ Neutre"

$ws1.Range("A5").Value = "BRVM - AUTRES SECTEURS"
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 2645.74
$ws1.Range("E5").Value = 657.73
$ws1.Range("F5").Value = "🟡 Observer"
$ws1.Range("G5").Value = "➖ Neutre"

$ws1.Range("A6").Value = "NEI-CEDA CI"
$ws1.Range("B6").Value = 0
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = 2370
$ws1.Range("E6").Value = 590
$ws1.Range("F6").Value = "🟡 Observer"
$ws1.Range("G6").Value = "➖ Neutre"

$ws1.Range("A7").Value = "SETAO CI"
$ws1.Range("B7").Value = 0
$ws1.Range("C7").Value = 4
$ws1.Range("D7").Value = 2360
$ws1.Range("E7").Value = 590
$ws1.Range("F7").Value = "🟡 Observer"
$ws1.Range("G7").Value = "➖ Neutre"

$ws1.Range("A8").Value = "UNIWAX CI"
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 4
$ws1.Range("D8").Value = 2290
$ws1.Range("E8").Value = 580
$ws1.Range("F8").Value = "🟡 Observer"
$ws1.Range("G8").Value = "➖ Neutre"

$ws1.Range("A9").Value = "AIR LIQUIDE CI"
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 4
$ws1.Range("D9").Value = 2105
$ws1.Range("E9").Value = 525
$ws1.Range("F9").Value = "🟡 Observer"
$ws1.Range("G9").Value = "➖ Neutre"

$ws1.Range("A10").Value = "BRVM - DISTRIBUTION"
$ws1.Range("B10").Value = 0
$ws1.Range("C10").Value = 4
$ws1.Range("D10").Value = 1430.32
$ws1.Range("E10").Value = 355.73
$ws1.Range("F10").Value = "🟡 Observer"
$ws1.Range("G10").Value = "➖ Neutre"

$ws1.Range("A11").Value = "BRVM - TRANSPORT"
$ws1.Range("B11").Value = 0
$ws1.Range("C11").Value = 4
$ws1.Range("D11").Value = 1387.86
$ws1.Range("E11").Value = 345.13
$ws1.Range("F11").Value = "🟡 Observer"
$ws1.Range("G11").Value = "➖ Neutre"

$ws1.Range("A12").Value = "BRVM - AGRICULTURE"
$ws1.Range("B12").Value = 0
$ws1.Range("C12").Value = 4
$ws1.Range("D12").Value = 1230.79
$ws1.Range("E12").Value = 306.08
$ws1.Range("F12").Value = "🟡 Observer"
$ws1.Range("G12").Value = "➖ Neutre"

$ws1.Range("A13").Value = "BRVM - INDUSTRIE"
$ws1.Range("B13").Value = 0
$ws1.Range("C13").Value = 4
$ws1.Range("D13").Value = 885.97
$ws1.Range("E13").Value = 231.32
$ws1.Range("F13").Value = "🟡 Observer"
$ws1.Range("G13").Value = "➖ Neutre"

$ws1.Range("A14").Value = "BRVM - CONSOMMATION DE BASE"
$ws1.Range("B14").Value = 0
$ws1.Range("C14").Value = 4
$ws1.Range("D14").Value = 745.73
$ws1.Range("E14").Value = 192.41
$ws1.Range("F14").Value = "🟡 Observer"
$ws1.Range("G14").Value = "➖ Neutre"

$ws1.Range("A15").Value = "BRVM-PRINCIPAL"
$ws1.Range("B15").Value = 0
$ws1.Range("C15").Value = 4
$ws1.Range("D15").Value = 720.77
$ws1.Range("E15").Value = 182.34
$ws1.Range("F15").Value = "🟡 Observer"
$ws1.Range("G15").Value = "➖ Neutre"

$ws1.Range("A16").Value = "BRVM - INDUSTRIELS"
$ws1.Range("B16").Value = 0
$ws1.Range("C16").Value = 4
$ws1.Range("D16").Value = 537.76
$ws1.Range("E16").Value = 138.34
$ws1.Range("F16").Value = "🟡 Observer"
$ws1.Range("G16").Value = "➖ Neutre"

$ws1.Range("A17").Value = "BRVM-PRESTIGE"
$ws1.Range("B17").Value = 0
$ws1.Range("C17").Value = 4
$ws1.Range("D17").Value = 518.31
$ws1.Range("E17").Value = 129.75
$ws1.Range("F17").Value = "🟡 Observer"
$ws1.Range("G17").Value = "➖ Neutre"

$ws1.Range("A18").Value = "BRVM - FINANCES"
$ws1.Range("B18").Value = 0
$ws1.Range("C18").Value = 4
$ws1.Range("D18").Value = 489.47
$ws1.Range("E18").Value = 123.03
$ws1.Range("F18").Value = "🟡 Observer"
$ws1.Range("G18").Value = "➖ Neutre"

$ws1.Range("A19").Value = "BRVM - SERVICES FINANCIERS"
$ws1.Range("B19").Value = 0
$ws1.Range("C19").Value = 4
$ws1.Range("D19").Value = 481.05
$ws1.Range("E19").Value = 120.91
$ws1.Range("F19").Value = "🟡 Observer"
$ws1.Range("G19").Value = "➖ Neutre"

$ws1.Range("A20").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Range("B20").Value = 0
$ws1.Range("C20").Value = 4
$ws1.Range("D20").Value = 425.13
$ws1.Range("E20").Value = 106.13
$ws1.Range("F20").Value = "🟡 Observer"
$ws1.Range("G20").Value = "➖ Neutre"

$ws1.Range("A21").Value = "BRVM - ENERGIE"
$ws1.Range("B21").Value = 0
$ws1.Range("C21").Value = 4
$ws1.Range("D21").Value = 420.42
$ws1.Range("E21").Value = 104.23
$ws1.Range("F21").Value = "🟡 Observer"
$ws1.Range("G21").Value = "➖ Neutre"

$ws1.Range("A22").Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Range("B22").Value = 0
$ws1.Range("C22").Value = 4
$ws1.Range("D22").Value = 374.21
$ws1.Range("E22").Value = 92.84
$ws1.Range("F22").Value = "🟡 Observer"
$ws1.Range("G22").Value = "➖ Neutre"

$ws1.Range("A23").Value = "UNILEVER CI (UNLC)"
$ws1.Range("B23").Value = 4
$ws1.Range("C23").Value = 0
$ws1.Range("D23").Value = 29.96
$ws1.Range("E23").Value = 7.49
$ws1.Range("F23").Value = "🟢 Achat"
$ws1.Range("G23").Value = "✅ Renforcer"

$ws1.Range("A24").Value = "FILTISAC CI (FTSC)"
$ws1.Range("B24").Value = 2
$ws1.Range("C24").Value = 0
$ws1.Range("D24").Value = 11.78
$ws1.Range("E24").Value = 7.45
$ws1.Range("F24").Value = "🟡 Observer"
$ws1.Range("G24").Value = "➖ Neutre"

$ws1.Range("A25").Value = "CIE CI (CIEC)"
$ws1.Range("B25").Value = 1
$ws1.Range("C25").Value = 0
$ws1.Range("D25").Value = 7.33
$ws1.Range("E25").Value = 7.33
$ws1.Range("F25").Value = "🟡 Observer"
$ws1.Range("G25").Value = "➖ Neutre"

$ws1.Range("A26").Value = "SODE CI (SDCC)"
$ws1.Range("B26").Value = 1
$ws1.Range("C26").Value = 0
$ws1.Range("D26").Value = 7.26
$ws1.Range("E26").Value = 7.26
$ws1.Range("F26").Value = "🟡 Observer"
$ws1.Range("G26").Value = "➖ Neutre"

$ws1.Range("A27").Value = "SITAB CI (STBC)"
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = 0
$ws1.Range("D27").Value = 6.35
$ws1.Range("E27").Value = 6.35
$ws1.Range("F27").Value = "🟡 Observer"
$ws1.Range("G27").Value = "➖ Neutre"

$ws1.Range("A28").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Range("B28").Value = 1
$ws1.Range("C28").Value = 0
$ws1.Range("D28").Value = 5.47
$ws1.Range("E28").Value = 5.47
$ws1.Range("F28").Value = "🟡 Observer"
$ws1.Range("G28").Value = "➖ Neutre"

$ws1.Range("A29").Value = "SMB CI (SMBC)"
$ws1.Range("B29").Value = 1
$ws1.Range("C29").Value = 0
$ws1.Range("D29").Value = 3.92
$ws1.Range("E29").Value = 3.92
$ws1.Range("F29").Value = "🟡 Observer"
$ws1.Range("G29").Value = "➖ Neutre"

$ws1.Range("A30").Value = "SOLIBRA CI (SLBC)"
$ws1.Range("B30").Value = 1
$ws1.Range("C30").Value = 0
$ws1.Range("D30").Value = 3.83
$ws1.Range("E30").Value = 3.83
$ws1.Range("F30").Value = "🟡 Observer"
$ws1.Range("G30").Value = "➖ Neutre"

$ws1.Range("A31").Value = "BERNABE CI (BNBC)"
$ws1.Range("B31").Value = 2
$ws1.Range("C31").Value = 1
$ws1.Range("D31").Value = 3.25
$ws1.Range("E31").Value = 7.2
$ws1.Range("F31").Value = "🟡 Observer"
$ws1.Range("G31").Value = "👀 À surveiller"

$ws1.Range("A32").Value = "BANK OF AFRICA NG (BOAN)"
$ws1.Range("B32").Value = 1
$ws1.Range("C32").Value = 1
$ws1.Range("D32").Value = 1.94
$ws1.Range("E32").Value = 5.1
$ws1.Range("F32").Value = "🟡 Observer"
$ws1.Range("G32").Value = "👀 À surveiller"

$ws1.Range("A33").Value = "CFAO MOTORS CI (CFAC)"
$ws1.Range("B33").Value = 1
$ws1.Range("C33").Value = 1
$ws1.Range("D33").Value = 1.63
$ws1.Range("E33").Value = 4.62
$ws1.Range("F33").Value = "🟡 Observer"
$ws1.Range("G33").Value = "👀 À surveiller"

$ws1.Range("A34").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Range("B34").Value = 1
$ws1.Range("C34").Value = 1
$ws1.Range("D34").Value = 0.37
$ws1.Range("E34").Value = -5.88
$ws1.Range("F34").Value = "🟡 Observer"
$ws1.Range("G34").Value = "👀 À surveiller"

$ws1.Range("A35").Value = "TOTAL"
$ws1.Range("B35").Value = 0
$ws1.Range("C35").Value = 4
$ws1.Range("D35").Value = 0
$ws1.Range("E35").Value = 0
$ws1.Range("F35").Value = "🟡 Observer"
$ws1.Range("G35").Value = "➖ Neutre"

$ws1.Range("A36").Value = "UNIWAX CI (UNXC)"
$ws1.Range("B36").Value = 1
$ws1.Range("C36").Value = 1
$ws1.Range("D36").Value = -0.71
$ws1.Range("E36").Value = -5.17
$ws1.Range("F36").Value = "🟡 Observer"
$ws1.Range("G36").Value = "👀 À surveiller"

$ws1.Range("A37").Value = "NESTLE CI (NTLC)"
$ws1.Range("B37").Value = 0
$ws1.Range("C37").Value = 1
$ws1.Range("D37").Value = -0.85
$ws1.Range("E37").Value = -0.85
$ws1.Range("F37").Value = "🟡 Observer"
$ws1.Range("G37").Value = "➖ Neutre"

$ws1.Range("A38").Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Range("B38").Value = 1
$ws1.Range("C38").Value = 1
$ws1.Range("D38").Value = -1.47
$ws1.Range("E38").Value = 4.46
$ws1.Range("F38").Value = "🟡 Observer"
$ws1.Range("G38").Value = "👀 À surveiller"

$ws1.Range("A39").Value = "NEI-CEDA CI (NEIC)"
$ws1.Range("B39").Value = 0
$ws1.Range("C39").Value = 1
$ws1.Range("D39").Value = -1.68
$ws1.Range("E39").Value = -1.68
$ws1.Range("F39").Value = "🟡 Observer"
$ws1.Range("G39").Value = "➖ Neutre"

$ws1.Range("A40").Value = "BANK OF AFRICA ML (BOAM)"
$ws1.Range("B40").Value = 0
$ws1.Range("C40").Value = 1
$ws1.Range("D40").Value = -1.91
$ws1.Range("E40").Value = -1.91
$ws1.Range("F40").Value = "🟡 Observer"
$ws1.Range("G40").Value = "➖ Neutre"

$ws1.Range("A41").Value = "VIVO ENERGY CI (SHEC)"
$ws1.Range("B41").Value = 0
$ws1.Range("C41").Value = 1
$ws1.Range("D41").Value = -2.4
$ws1.Range("E41").Value = -2.4
$ws1.Range("F41").Value = "🟡 Observer"
$ws1.Range("G41").Value = "➖ Neutre"

$ws1.Range("A42").Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Range("B42").Value = 0
$ws1.Range("C42").Value = 1
$ws1.Range("D42").Value = -2.42
$ws1.Range("E42").Value = -2.42
$ws1.Range("F42").Value = "🟡 Observer"
$ws1.Range("G42").Value = "➖ Neutre"

$ws1.Range("A43").Value = "AIR LIQUIDE CI (SIVC)"
$ws1.Range("B43").Value = 0
$ws1.Range("C43").Value = 1
$ws1.Range("D43").Value = -2.86
$ws1.Range("E43").Value = -2.86
$ws1.Range("F43").Value = "🟡 Observer"
$ws1.Range("G43").Value = "➖ Neutre"

$ws1.Range("A44").Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Range("B44").Value = 0
$ws1.Range("C44").Value = 1
$ws1.Range("D44").Value = -4.03
$ws1.Range("E44").Value = -4.03
$ws1.Range("F44").Value = "🟡 Observer"
$ws1.Range("G44").Value = "➖ Neutre"

$ws1.Range("A45").Value = "ORANGE COTE D'IVOIRE (ORAC)"
$ws1.Range("B45").Value = 0
$ws1.Range("C45").Value = 2
$ws1.Range("D45").Value = -4.09
$ws1.Range("E45").Value = -2.02
$ws1.Range("F45").Value = "🟡 Observer"
$ws1.Range("G45").Value = "➖ Neutre"

$ws1.Range("A46").Value = "SETAO CI (STAC)"
$ws1.Range("B46").Value = 0
$ws1.Range("C46").Value = 1
$ws1.Range("D46").Value = -4.2
$ws1.Range("E46").Value = -4.2
$ws1.Range("F46").Value = "🟡 Observer"
$ws1.Range("G46").Value = "➖ Neutre"

$ws1.Range("A47").Value = "ONATEL BF (ONTBF)"
$ws1.Range("B47").Value = 1
$ws1.Range("C47").Value = 2
$ws1.Range("D47").Value = -4.65
$ws1.Range("E47").Value = -4.37
$ws1.Range("F47").Value = "🟡 Observer"
$ws1.Range("G47").Value = "👀 À surveiller"

$ws1.Range("A48").Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws1.Range("B48").Value = 0
$ws1.Range("C48").Value = 2
$ws1.Range("D48").Value = -7.04
$ws1.Range("E48").Value = -4.47
$ws1.Range("F48").Value = "🟡 Observer"
$ws1.Range("G48").Value = "➖ Neutre"

# --- Sheet2: Top_YTD (rows 2-11, column B only) ---
$ws2.Range("B2").Value = 8830336.91
$ws2.Range("B3").Value = 360605.15
$ws2.Range("B4").Value = 355799.37
$ws2.Range("B5").Value = 336042.96
$ws2.Range("B6").Value = 229856.3
$ws2.Range("B7").Value = 226559.31
$ws2.Range("B8").Value = 204339.62
$ws2.Range("B9").Value = 153698.75
$ws2.Range("B10").Value = 43734.31
$ws2.Range("B11").Value = 39809.79
